{"js": "// PUNTO 1 - TP3: fix the pseudocode for buscarElemento()\n//  1) \"ind=ind+1\"   -> \"ind=ind+1;\"   (missing trailing semicolon)\n//  2) final \"retorna ind;\" -> \"retorna -1;\" (function must return -1, not ind,\n//     when the element was not found) and a blank line is left after it.\n\n// --- 1) Append the missing semicolon after \"ind=ind+1\" ---------------------\nconst incResults = context.document.body.search(\"ind=ind+1\", { matchCase: true, matchWildcards: false });\nincResults.load(\"items\");\nawait context.sync();\n\nif (incResults.items.length > 0) {\n  // Insert right at the end of the matched text so it lands in the same run\n  // (inheriting its Times New Roman formatting) instead of creating a new,\n  // unformatted trailing run.\n  incResults.items[0].insertText(\";\", Word.InsertLocation.end);\n}\n\nawait context.sync();\n\n// --- 2) Replace the trailing \"retorna ind;\" with \"retorna -1;\" -------------\nconst retResults = context.document.body.search(\"retorna ind;\", { matchCase: true, matchWildcards: false });\nretResults.load(\"items\");\nawait context.sync();\n\nif (retResults.items.length > 0) {\n  // There are two \"retorna ind;\" occurrences in the document; the one that\n  // needs to change is the very last one (the function's final return).\n  const last = retResults.items[retResults.items.length - 1];\n  last.insertText(\"retorna -1;\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Insert a new, empty paragraph right after that line (matching the\n  // formatting of the surrounding text) leaving a blank line before the\n  // document's final closing paragraph.\n  const ownerParagraph = last.paragraphs.getFirst();\n  ownerParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# PUNTO 1 - TP3: fix the pseudocode for buscarElemento()\n#  1) \"ind=ind+1\"   -> \"ind=ind+1;\"   (missing trailing semicolon)\n#  2) final \"retorna ind;\" -> \"retorna -1;\" (function must return -1, not ind,\n#     when the element was not found) and a blank line is left after it.\n\n$d = $word.ActiveDocument\n\n# --- 1) Append the missing semicolon after \"ind=ind+1\" ----------------------\n$incRange = $d.Content\n$incRange.Find.ClearFormatting()\nif ($incRange.Find.Execute(\"ind=ind+1\")) {\n    # Re-writing Text on the matched range keeps it inside the existing\n    # run/formatting (Times New Roman) instead of appending an unformatted\n    # trailing run.\n    $incRange.Text = \"ind=ind+1;\"\n}\n\n# --- 2) Replace the trailing \"retorna ind;\" with \"retorna -1;\" -------------\n# There are two \"retorna ind;\" occurrences in the document; the one that\n# needs to change is the very last one (the function's final return), so we\n# walk every match and remember the last span found.\n$scanRange = $d.Content\n$scanRange.Find.ClearFormatting()\n$lastStart = -1\n$lastEnd = -1\nwhile ($scanRange.Find.Execute(\"retorna ind;\")) {\n    $lastStart = $scanRange.Start\n    $lastEnd = $scanRange.End\n    $scanRange.Collapse(0)\n}\n\nif ($lastStart -ge 0) {\n    $target = $d.Range($lastStart, $lastEnd)\n    $target.Text = \"retorna -1;\"\n\n    # Insert a new, empty paragraph right after that line (matching the\n    # formatting of the surrounding text) leaving a blank line before the\n    # document's final closing paragraph.\n    $target.InsertParagraphAfter()\n}\n"}
